{"js": "// Replace the 25 division-problem values in the table with the new values,\n// matching the unified diff exactly. Each \"from\" text is unique in the\n// document, so a targeted search + full-text replace for each pair is safe.\nconst replacements = [\n  [\"362\u00f76=\", \"179\u00f75=\"],\n  [\"263\u00f72=\", \"137\u00f72=\"],\n  [\"690\u00f76=\", \"856\u00f79=\"],\n  [\"939\u00f79=\", \"161\u00f78=\"],\n  [\"741\u00f74=\", \"519\u00f78=\"],\n  [\"484\u00f78=\", \"106\u00f72=\"],\n  [\"844\u00f78=\", \"978\u00f72=\"],\n  [\"659\u00f75=\", \"187\u00f72=\"],\n  [\"812\u00f79=\", \"487\u00f74=\"],\n  [\"810\u00f79=\", \"545\u00f75=\"],\n  [\"984\u00f72=\", \"531\u00f72=\"],\n  [\"643\u00f75=\", \"558\u00f79=\"],\n  [\"830\u00f78=\", \"510\u00f73=\"],\n  [\"456\u00f76=\", \"414\u00f79=\"],\n  [\"180\u00f74=\", \"439\u00f74=\"],\n  [\"629\u00f77=\", \"421\u00f73=\"],\n  [\"541\u00f72=\", \"331\u00f79=\"],\n  [\"875\u00f79=\", \"154\u00f79=\"],\n  [\"577\u00f75=\", \"154\u00f75=\"],\n  [\"965\u00f79=\", \"645\u00f76=\"],\n  [\"235\u00f79=\", \"495\u00f78=\"],\n  [\"904\u00f77=\", \"916\u00f74=\"],\n  [\"532\u00f74=\", \"153\u00f78=\"],\n  [\"204\u00f77=\", \"461\u00f79=\"],\n  [\"275\u00f74=\", \"876\u00f75=\"],\n];\n\nfor (const [oldText, newText] of replacements) {\n  const results = context.document.body.search(oldText, { matchCase: true, matchWholeWord: false });\n  results.load(\"items\");\n  await context.sync();\n\n  for (const item of results.items) {\n    item.insertText(newText, \"Replace\");\n  }\n  await context.sync();\n}\n", "ps1": "# Replace the 25 division-problem values in the table with the new values,\n# matching the unified diff exactly. Each \"from\" text is unique in the\n# document, so a targeted Find/Replace (ReplaceAll) for each pair is safe.\n$d = $word.ActiveDocument\n\n$replacements = @(\n    @(\"362\u00f76=\", \"179\u00f75=\"),\n    @(\"263\u00f72=\", \"137\u00f72=\"),\n    @(\"690\u00f76=\", \"856\u00f79=\"),\n    @(\"939\u00f79=\", \"161\u00f78=\"),\n    @(\"741\u00f74=\", \"519\u00f78=\"),\n    @(\"484\u00f78=\", \"106\u00f72=\"),\n    @(\"844\u00f78=\", \"978\u00f72=\"),\n    @(\"659\u00f75=\", \"187\u00f72=\"),\n    @(\"812\u00f79=\", \"487\u00f74=\"),\n    @(\"810\u00f79=\", \"545\u00f75=\"),\n    @(\"984\u00f72=\", \"531\u00f72=\"),\n    @(\"643\u00f75=\", \"558\u00f79=\"),\n    @(\"830\u00f78=\", \"510\u00f73=\"),\n    @(\"456\u00f76=\", \"414\u00f79=\"),\n    @(\"180\u00f74=\", \"439\u00f74=\"),\n    @(\"629\u00f77=\", \"421\u00f73=\"),\n    @(\"541\u00f72=\", \"331\u00f79=\"),\n    @(\"875\u00f79=\", \"154\u00f79=\"),\n    @(\"577\u00f75=\", \"154\u00f75=\"),\n    @(\"965\u00f79=\", \"645\u00f76=\"),\n    @(\"235\u00f79=\", \"495\u00f78=\"),\n    @(\"904\u00f77=\", \"916\u00f74=\"),\n    @(\"532\u00f74=\", \"153\u00f78=\"),\n    @(\"204\u00f77=\", \"461\u00f79=\"),\n    @(\"275\u00f74=\", \"876\u00f75=\")\n)\n\nforeach ($pair in $replacements) {\n    $oldText = $pair[0]\n    $newText = $pair[1]\n\n    $find = $d.Content.Find\n    $find.ClearFormatting()\n    $find.Replacement.ClearFormatting()\n    $find.Text = $oldText\n    $find.Replacement.Text = $newText\n    $find.Forward = $true\n    $find.Wrap = 1  # wdFindContinue\n    $find.Execute($find.Text, $false, $false, $false, $false, $false, $true, 1, $false, $find.Replacement.Text, 2)  # wdReplaceAll\n}\n"}
